$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header/summary fields ---
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:48 AM"
$ws.Range("C8").Value = 18871.41
$ws.Range("C9").Value = 51
$ws.Range("C10").Value = "07/28/2025 to 08/03/25"

# Re-assert an existing column-I merge so the sheet dimension/used-range
# keeps extending through column I (matches the source report layout).
$ws.Range("G8:I8").Merge()

# --- Helper: copy formatting only (xlPasteFormats = -4122) ---
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}

# --- Block: Wednesday (07/30/2025) ---
Copy-Format "A14" "A49"
$ws.Cells.Item(49,1).Value = "Wednesday (07/30/2025)"
$ws.Range("A49:H49").Merge()

Copy-Format "A15:H15" "A50:H50"
$ws.Cells.Item(50,1).Value = "Point Number"
$ws.Cells.Item(50,2).Value = "Billable Unit Code"
$ws.Cells.Item(50,3).Value = "Work Type"
$ws.Cells.Item(50,4).Value = "Unit Description"
$ws.Cells.Item(50,5).Value = "Unit of Measure"
$ws.Cells.Item(50,6).Value = "# Units"
$ws.Cells.Item(50,7).Value = "N/A"
$ws.Cells.Item(50,8).Value = "Pricing"

Copy-Format "A16:H16" "A51:H51"
$ws.Cells.Item(51,1).Value = "Point 19"
$ws.Cells.Item(51,2).Value = "ANC-DHM-10-84-T1-C"
$ws.Cells.Item(51,3).Value = "Inst"
$ws.Cells.Item(51,4).Value = "ANC,Dbl Hlx Mach,10in,84in,TpEye 1in,Cor"
$ws.Cells.Item(51,5).Value = "EA"
$ws.Cells.Item(51,6).Value = 2
$ws.Cells.Item(51,8).Value = 435.06

Copy-Format "A17:H17" "A52:H52"
$ws.Cells.Item(52,1).Value = "Point 19"
$ws.Cells.Item(52,2).Value = "BKT-IP8-F-C"
$ws.Cells.Item(52,3).Value = "Inst"
$ws.Cells.Item(52,4).Value = "BKT,Insulator Post 8in,Fbrgls,Corrosive"
$ws.Cells.Item(52,5).Value = "EA"
$ws.Cells.Item(52,6).Value = 3
$ws.Cells.Item(52,8).Value = 95.16

Copy-Format "A16:H16" "A53:H53"
$ws.Cells.Item(53,1).Value = "Point 19"
$ws.Cells.Item(53,2).Value = "CNC-SNB-40"
$ws.Cells.Item(53,3).Value = "Inst"
$ws.Cells.Item(53,4).Value = "CNC,Splice Non-Tension Bare,#1/0-#4/0"
$ws.Cells.Item(53,5).Value = "EA"
$ws.Cells.Item(53,6).Value = 5
$ws.Cells.Item(53,8).Value = 101.45

Copy-Format "A17:H17" "A54:H54"
$ws.Cells.Item(54,1).Value = "Point 19"
$ws.Cells.Item(54,2).Value = "DEC-20AL-C"
$ws.Cells.Item(54,3).Value = "Inst"
$ws.Cells.Item(54,4).Value = "DEC,#4 - #2/0 AA,AL,AS,Corrosive"
$ws.Cells.Item(54,5).Value = "EA"
$ws.Cells.Item(54,6).Value = 9
$ws.Cells.Item(54,8).Value = 2576.25

Copy-Format "A16:H16" "A55:H55"
$ws.Cells.Item(55,1).Value = "Point 19"
$ws.Cells.Item(55,2).Value = "GYF-38-D-78P-EP-C"
$ws.Cells.Item(55,3).Value = "Inst"
$ws.Cells.Item(55,4).Value = "GYF,3/8,Down,78in Pole mt,EyePlate,Corr"
$ws.Cells.Item(55,5).Value = "EA"
$ws.Cells.Item(55,6).Value = 3
$ws.Cells.Item(55,8).Value = 238.05

Copy-Format "A17:H17" "A56:H56"
$ws.Cells.Item(56,1).Value = "Point 19"
$ws.Cells.Item(56,2).Value = "INS-15-P-S-C"
$ws.Cells.Item(56,3).Value = "Inst"
$ws.Cells.Item(56,4).Value = "INS,15kV,Pin,Silicon Polymer,Corr"
$ws.Cells.Item(56,5).Value = "EA"
$ws.Cells.Item(56,6).Value = 3
$ws.Cells.Item(56,8).Value = 282.51

Copy-Format "A16:H16" "A57:H57"
$ws.Cells.Item(57,1).Value = "Point 19"
$ws.Cells.Item(57,2).Value = "POL-45-2"
$ws.Cells.Item(57,3).Value = "Inst"
$ws.Cells.Item(57,4).Value = "Pole,45ft,Class 2"
$ws.Cells.Item(57,5).Value = "EA"
$ws.Cells.Item(57,6).Value = 1
$ws.Cells.Item(57,8).Value = 478.55

Copy-Format "A17:H17" "A58:H58"
$ws.Cells.Item(58,1).Value = "Point 19"
$ws.Cells.Item(58,2).Value = "SAA-DE-20-C"
$ws.Cells.Item(58,3).Value = "Inst"
$ws.Cells.Item(58,4).Value = "SAA,DE Clamp #4-2/0, Corr."
$ws.Cells.Item(58,5).Value = "EA"
$ws.Cells.Item(58,6).Value = 3
$ws.Cells.Item(58,8).Value = 165.54

Copy-Format "A16:H16" "A59:H59"
$ws.Cells.Item(59,1).Value = "Point 21"
$ws.Cells.Item(59,2).Value = "ANC-DHM-10-84-T1-C"
$ws.Cells.Item(59,3).Value = "Inst"
$ws.Cells.Item(59,4).Value = "ANC,Dbl Hlx Mach,10in,84in,TpEye 1in,Cor"
$ws.Cells.Item(59,5).Value = "EA"
$ws.Cells.Item(59,6).Value = 2
$ws.Cells.Item(59,8).Value = 435.06

Copy-Format "A27" "A60"
Copy-Format "H27" "H60"
$ws.Cells.Item(60,1).Value = "TOTAL"
$ws.Cells.Item(60,8).Value = 4807.630000000001
$ws.Range("A60:G60").Merge()

# --- Block: Thursday (07/31/2025) ---
Copy-Format "A30" "A63"
$ws.Cells.Item(63,1).Value = "Thursday (07/31/2025)"
$ws.Range("A63:H63").Merge()

Copy-Format "A31:H31" "A64:H64"
$ws.Cells.Item(64,1).Value = "Point Number"
$ws.Cells.Item(64,2).Value = "Billable Unit Code"
$ws.Cells.Item(64,3).Value = "Work Type"
$ws.Cells.Item(64,4).Value = "Unit Description"
$ws.Cells.Item(64,5).Value = "Unit of Measure"
$ws.Cells.Item(64,6).Value = "# Units"
$ws.Cells.Item(64,7).Value = "N/A"
$ws.Cells.Item(64,8).Value = "Pricing"

Copy-Format "A32:H32" "A65:H65"
$ws.Cells.Item(65,1).Value = "Point 28"
$ws.Cells.Item(65,2).Value = "CON-2-AAA-1-B-REEL"
$ws.Cells.Item(65,3).Value = "Rem"
$ws.Cells.Item(65,4).Value = "CON,#2 AWG,Alum Alloy,One,Bare,Reels"
$ws.Cells.Item(65,5).Value = "FT"
$ws.Cells.Item(65,6).Value = 464
$ws.Cells.Item(65,8).Value = 431.52

Copy-Format "A33:H33" "A66:H66"
$ws.Cells.Item(66,1).Value = "Point 28"
$ws.Cells.Item(66,2).Value = "EQL-1-4-C-2-S-X-C"
$ws.Cells.Item(66,3).Value = "Rem"
$ws.Cells.Item(66,4).Value = "EQL,1 Ph,#4,CU Sol,#2,CU Str,Xfr,Corr"
$ws.Cells.Item(66,5).Value = "EA"
$ws.Cells.Item(66,6).Value = 1
$ws.Cells.Item(66,8).Value = 116

Copy-Format "A32:H32" "A67:H67"
$ws.Cells.Item(67,1).Value = "Point 28"
$ws.Cells.Item(67,2).Value = "POL-40-4"
$ws.Cells.Item(67,3).Value = "Rem"
$ws.Cells.Item(67,4).Value = "Pole,40ft,Class 4"
$ws.Cells.Item(67,5).Value = "EA"
$ws.Cells.Item(67,6).Value = 1
$ws.Cells.Item(67,8).Value = 198.88

Copy-Format "A33:H33" "A68:H68"
$ws.Cells.Item(68,1).Value = "Point 30"
$ws.Cells.Item(68,2).Value = "CON-2-AAA-1-B-REEL"
$ws.Cells.Item(68,3).Value = "Rem"
$ws.Cells.Item(68,4).Value = "CON,#2 AWG,Alum Alloy,One,Bare,Reels"
$ws.Cells.Item(68,5).Value = "FT"
$ws.Cells.Item(68,6).Value = 1312
$ws.Cells.Item(68,8).Value = 1220.16

Copy-Format "A32:H32" "A69:H69"
$ws.Cells.Item(69,1).Value = "Point 32"
$ws.Cells.Item(69,2).Value = "CON-2-AAA-1-B-REEL"
$ws.Cells.Item(69,3).Value = "Rem"
$ws.Cells.Item(69,4).Value = "CON,#2 AWG,Alum Alloy,One,Bare,Reels"
$ws.Cells.Item(69,5).Value = "FT"
$ws.Cells.Item(69,6).Value = 880
$ws.Cells.Item(69,8).Value = 818.4

Copy-Format "A33:H33" "A70:H70"
$ws.Cells.Item(70,1).Value = "Point 01"
$ws.Cells.Item(70,2).Value = "PLA-DLOC"
$ws.Cells.Item(70,3).Value = "Inst"
$ws.Cells.Item(70,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(70,5).Value = "EA"
$ws.Cells.Item(70,6).Value = 4
$ws.Cells.Item(70,8).Value = 476.4

Copy-Format "A32:H32" "A71:H71"
$ws.Cells.Item(71,1).Value = "Point 01"
$ws.Cells.Item(71,2).Value = "PLA-HDIG"
$ws.Cells.Item(71,3).Value = "Inst"
$ws.Cells.Item(71,4).Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Cells.Item(71,5).Value = "EA"
$ws.Cells.Item(71,6).Value = 1
$ws.Cells.Item(71,8).Value = 648.53

Copy-Format "A33:H33" "A72:H72"
$ws.Cells.Item(72,1).Value = "Point 03"
$ws.Cells.Item(72,2).Value = "PLA-HDIG"
$ws.Cells.Item(72,3).Value = "Inst"
$ws.Cells.Item(72,4).Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Cells.Item(72,5).Value = "EA"
$ws.Cells.Item(72,6).Value = 1
$ws.Cells.Item(72,8).Value = 648.53

Copy-Format "A32:H32" "A73:H73"
$ws.Cells.Item(73,1).Value = "Point 19"
$ws.Cells.Item(73,2).Value = "PLA-DLOC"
$ws.Cells.Item(73,3).Value = "Inst"
$ws.Cells.Item(73,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(73,5).Value = "EA"
$ws.Cells.Item(73,6).Value = 2
$ws.Cells.Item(73,8).Value = 238.2

Copy-Format "A33:H33" "A74:H74"
$ws.Cells.Item(74,1).Value = "Point 23"
$ws.Cells.Item(74,2).Value = "PLA-DLOC"
$ws.Cells.Item(74,3).Value = "Inst"
$ws.Cells.Item(74,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(74,5).Value = "EA"
$ws.Cells.Item(74,6).Value = 4
$ws.Cells.Item(74,8).Value = 476.4

Copy-Format "A32:H32" "A75:H75"
$ws.Cells.Item(75,1).Value = "Point 25"
$ws.Cells.Item(75,2).Value = "PLA-DLOC"
$ws.Cells.Item(75,3).Value = "Inst"
$ws.Cells.Item(75,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(75,5).Value = "EA"
$ws.Cells.Item(75,6).Value = 2
$ws.Cells.Item(75,8).Value = 238.2

Copy-Format "A33:H33" "A76:H76"
$ws.Cells.Item(76,1).Value = "Point 27"
$ws.Cells.Item(76,2).Value = "PLA-DLOC"
$ws.Cells.Item(76,3).Value = "Inst"
$ws.Cells.Item(76,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(76,5).Value = "EA"
$ws.Cells.Item(76,6).Value = 2
$ws.Cells.Item(76,8).Value = 238.2

Copy-Format "A32:H32" "A77:H77"
$ws.Cells.Item(77,1).Value = "Point 31"
$ws.Cells.Item(77,2).Value = "PLA-DLOC"
$ws.Cells.Item(77,3).Value = "Inst"
$ws.Cells.Item(77,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(77,5).Value = "EA"
$ws.Cells.Item(77,6).Value = 4
$ws.Cells.Item(77,8).Value = 476.4

Copy-Format "A33:H33" "A78:H78"
$ws.Cells.Item(78,1).Value = "Point 19"
$ws.Cells.Item(78,2).Value = "PLA-BACK"
$ws.Cells.Item(78,3).Value = "Inst"
$ws.Cells.Item(78,4).Value = "Difficult Location Equip Adder-Backyard"
$ws.Cells.Item(78,5).Value = "EA"
$ws.Cells.Item(78,6).Value = 4
$ws.Cells.Item(78,8).Value = 476.4

Copy-Format "A32:H32" "A79:H79"
$ws.Cells.Item(79,1).Value = "Point 23"
$ws.Cells.Item(79,2).Value = "PLA-BACK"
$ws.Cells.Item(79,3).Value = "Inst"
$ws.Cells.Item(79,4).Value = "Difficult Location Equip Adder-Backyard"
$ws.Cells.Item(79,5).Value = "EA"
$ws.Cells.Item(79,6).Value = 4
$ws.Cells.Item(79,8).Value = 476.4

Copy-Format "A33:H33" "A80:H80"
$ws.Cells.Item(80,1).Value = "Point 03"
$ws.Cells.Item(80,2).Value = "PLA-DLOC"
$ws.Cells.Item(80,3).Value = "Inst"
$ws.Cells.Item(80,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(80,5).Value = "EA"
$ws.Cells.Item(80,6).Value = 4
$ws.Cells.Item(80,8).Value = 476.4

Copy-Format "A32:H32" "A81:H81"
$ws.Cells.Item(81,1).Value = "Point 28"
$ws.Cells.Item(81,2).Value = "PLA-DLOC"
$ws.Cells.Item(81,3).Value = "Inst"
$ws.Cells.Item(81,4).Value = "PLA,Difficult Location"
$ws.Cells.Item(81,5).Value = "EA"
$ws.Cells.Item(81,6).Value = 4
$ws.Cells.Item(81,8).Value = 476.4

Copy-Format "A46" "A82"
Copy-Format "H46" "H82"
$ws.Cells.Item(82,1).Value = "TOTAL"
$ws.Cells.Item(82,8).Value = 8131.419999999997
$ws.Range("A82:G82").Merge()
